$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.116.92"
$ws.Range("D3").Value = "1.907.68"
$ws.Range("E3").Value = "  +5.24%  "
$ws.Range("D4").Formula = "'0.9993"
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Formula = "'252.63"
$ws.Range("E5").Value = "  +1.85%  "
$ws.Range("D6").Formula = "'0.9991"
$ws.Range("E6").Value = "  -0.01%  "
$ws.Range("D7").Formula = "'0.5091"
$ws.Range("E7").Value = "  +2.86%  "
$ws.Range("D8").Formula = "'45.25"
$ws.Range("E8").Value = "  +4.52%  "
$ws.Range("D9").Formula = "'0.3011"
$ws.Range("E9").Value = "  +8.10%  "
$ws.Range("D10").Formula = "'0.06841"
$ws.Range("E10").Value = "  +6.46%  "
$ws.Range("D11").Value = "1.906.15"
$ws.Range("E11").Value = "  +5.18%  "
$ws.Range("E12").Value = "  +2.81%  "
$ws.Range("D13").Formula = "'0.07324"
$ws.Range("E13").Value = "  +3.51%  "
$ws.Range("D14").Formula = "'0.6918"
$ws.Range("D15").Formula = "'87.00"
$ws.Range("E15").Value = "  +3.60%  "
$ws.Range("D16").Formula = "'4.922"
$ws.Range("E16").Value = "  +4.91%  "
$ws.Range("D17").Formula = "'0.000008363"
$ws.Range("E17").Value = "  +13.62%  "
$ws.Range("D18").Value = "30.117.71"
$ws.Range("E18").Value = "  +4.12%  "
$ws.Range("D19").Formula = "'0.9991"
$ws.Range("E19").Value = "  +0.00%  "
$ws.Range("E20").Value = "  +6.26%  "
$ws.Range("D21").Value = "2.151.23"
$ws.Range("E21").Value = "  +5.34%  "
$ws.Range("D22").Formula = "'0.9984"
$ws.Range("E22").Value = "  -0.06%  "
$ws.Range("D23").Formula = "'4.821"
$ws.Range("E23").Value = "  +5.13%  "
$ws.Range("D24").Formula = "'5.752"
$ws.Range("E24").Value = "  +7.23%  "
$ws.Range("D25").Formula = "'9.302"
$ws.Range("E25").Value = "  +5.25%  "
$ws.Range("D26").Formula = "'147.54"
$ws.Range("E26").Value = "  +3.36%  "
$ws.Range("D27").Formula = "'134.97"
$ws.Range("E27").Value = "  +4.20%  "
$ws.Range("E28").Value = "  +4.38%  "
$ws.Range("E29").Value = "  +5.77%  "
$ws.Range("D30").Formula = "'1.402"
$ws.Range("E30").Value = "  -0.92%  "
$ws.Range("D31").Formula = "'4.283"
$ws.Range("E31").Value = "  +3.24%  "
$ws.Range("D32").Formula = "'0.08857"
$ws.Range("E32").Value = "  +5.93%  "
$ws.Range("D33").Formula = "'4.003"
$ws.Range("E33").Value = "  +4.85%  "
$ws.Range("D34").Formula = "'0.05061"
$ws.Range("E34").Value = "  +2.03%  "
$ws.Range("D35").Formula = "'1.145"
$ws.Range("E35").Value = "  +3.43%  "
$ws.Range("D36").Formula = "'0.7233"
$ws.Range("E36").Value = "  +7.33%  "
$ws.Range("D37").Formula = "'2.690"
$ws.Range("E37").Value = "  +0.63%  "
$ws.Range("D38").Formula = "'2.824"
$ws.Range("E38").Value = "  +2.61%  "
$ws.Range("D39").Formula = "'2.277"
$ws.Range("E39").Value = "  -2.09%  "
$ws.Range("D40").Formula = "'0.9632"
$ws.Range("E40").Value = "  +1.05%  "
$ws.Range("D41").Formula = "'0.01691"
$ws.Range("E41").Value = "  +6.03%  "
$ws.Range("D42").Formula = "'6.100"
$ws.Range("E42").Value = "  -1.10%  "
$ws.Range("D43").Formula = "'0.4311"
$ws.Range("E43").Value = "  +5.04%  "
$ws.Range("D44").Formula = "'104.78"
$ws.Range("E44").Value = "  +4.38%  "
$ws.Range("D45").Formula = "'0.9990"
$ws.Range("E45").Value = "  +0.00%  "
$ws.Range("D46").Formula = "'7.678"
$ws.Range("E46").Value = "  +7.11%  "
$ws.Range("D47").Formula = "'0.1282"
$ws.Range("E47").Value = "  +4.83%  "
$ws.Range("D48").Formula = "'0.05754"
$ws.Range("E48").Value = "  +4.22%  "
$ws.Range("D49").Formula = "'33.27"
$ws.Range("E49").Value = "  +4.69%  "
$ws.Range("D50").Formula = "'8.435"
$ws.Range("E50").Value = "  +3.75%  "
$ws.Range("D51").Formula = "'0.3821"
$ws.Range("E51").Value = "  +5.00%  "
